$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 74, shifting rows 74-168 down to 75-169.
$ws.Rows.Item(74).Insert()

# Populate the newly inserted row 74 with the new weekly data point.
$ws.Cells.Item(74, 1).Value = 7
$ws.Cells.Item(74, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(74, 3).Value = "Ñuble"
$ws.Cells.Item(74, 4).Value = 45195
$ws.Cells.Item(74, 5).Value = 16
$ws.Cells.Item(74, 6).Value = 100112031
$ws.Cells.Item(74, 7).Value = "Poroto verde"
$ws.Cells.Item(74, 8).Value = "Magnum"
$ws.Cells.Item(74, 9).Value = "Primera"
$ws.Cells.Item(74, 10).Value = 50
$ws.Cells.Item(74, 11).Value = 30000
$ws.Cells.Item(74, 12).Value = 30000
$ws.Cells.Item(74, 13).Value = 30000
$ws.Cells.Item(74, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(74, 15).Value = "Perú"
$ws.Cells.Item(74, 16).Value = 1200
$ws.Cells.Item(74, 17).Value = 25
$ws.Cells.Item(74, 18).Value = "Hortaliza"
